$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting the existing data (rows 5-41) down to (rows 6-42)
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new weekly record
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44473
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 100112052
$ws.Range("G5").Value = "Albahaca"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 3500
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 3750
$ws.Range("N5").Value = "$/paquete"
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 3750
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
